$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the first data row (ID "H 72"); remaining rows shift up by one
$ws.Rows.Item(2).Delete()

# Overwrite rows 2-62 with the updated dataset / missing-value mask
$ws.Range("A2").Value = "H 1968 Sk 4"
$ws.Range("B2").Value = -20
$ws.Range("C2").Value = 10
$ws.Range("D2").Value = -14.8
$ws.Range("E2").Value = ""
$ws.Range("F2").Value = ""

$ws.Range("A3").Value = "H 12640 K XIII/3"
$ws.Range("B3").Value = -19.8
$ws.Range("C3").Value = 11.7
$ws.Range("D3").Value = -13.5
$ws.Range("E3").Value = -7.5
$ws.Range("F3").Value = 0.70917

$ws.Range("A4").Value = "RM 2"
$ws.Range("B4").Value = -19.7
$ws.Range("C4").Value = ""
$ws.Range("D4").Value = -13.5
$ws.Range("E4").Value = -7.2
$ws.Range("F4").Value = 0.70909

$ws.Range("A5").Value = "RM 8"
$ws.Range("B5").Value = -19.7
$ws.Range("C5").Value = 11.2
$ws.Range("D5").Value = ""
$ws.Range("E5").Value = -5.7
$ws.Range("F5").Value = 0.71055

$ws.Range("A6").Value = "RM 9"
$ws.Range("B6").Value = -18.7
$ws.Range("C6").Value = ""
$ws.Range("D6").Value = -15.4
$ws.Range("E6").Value = -6.4
$ws.Range("F6").Value = 0.71067

$ws.Range("A7").Value = "RM 14"
$ws.Range("B7").Value = -19.5
$ws.Range("C7").Value = ""
$ws.Range("D7").Value = -14.4
$ws.Range("E7").Value = -5
$ws.Range("F7").Value = 0.71266

$ws.Range("A8").Value = "RM 21"
$ws.Range("B8").Value = -19.8
$ws.Range("C8").Value = 15.1
$ws.Range("D8").Value = -14.2
$ws.Range("E8").Value = ""
$ws.Range("F8").Value = 0.70993

$ws.Range("A9").Value = "RM 28"
$ws.Range("B9").Value = -19.9
$ws.Range("C9").Value = 12
$ws.Range("D9").Value = -15.4
$ws.Range("E9").Value = -6.2
$ws.Range("F9").Value = 0.71194

$ws.Range("A10").Value = "RM 32"
$ws.Range("B10").Value = -19.5
$ws.Range("C10").Value = 15
$ws.Range("D10").Value = -13.8
$ws.Range("E10").Value = -7.1
$ws.Range("F10").Value = ""

$ws.Range("A11").Value = "RM 38"
$ws.Range("B11").Value = -19.9
$ws.Range("C11").Value = 15.5
$ws.Range("D11").Value = -13.9
$ws.Range("E11").Value = -6.6
$ws.Range("F11").Value = 0.71122

$ws.Range("A12").Value = "RM 42"
$ws.Range("B12").Value = -20.6
$ws.Range("C12").Value = 10.5
$ws.Range("D12").Value = ""
$ws.Range("E12").Value = -6.8
$ws.Range("F12").Value = 0.71118

$ws.Range("A13").Value = "RM 52 a"
$ws.Range("B13").Value = -19.8
$ws.Range("C13").Value = 11.5
$ws.Range("D13").Value = ""
$ws.Range("E13").Value = ""
$ws.Range("F13").Value = 0.71103

$ws.Range("A14").Value = "RM 58"
$ws.Range("B14").Value = -20.8
$ws.Range("C14").Value = 11.4
$ws.Range("D14").Value = -15.5
$ws.Range("E14").Value = -7.9
$ws.Range("F14").Value = 0.70989

$ws.Range("A15").Value = "RM 78"
$ws.Range("B15").Value = -18.9
$ws.Range("C15").Value = 12
$ws.Range("D15").Value = -15.4
$ws.Range("E15").Value = -12
$ws.Range("F15").Value = ""

$ws.Range("A16").Value = "RM 81"
$ws.Range("B16").Value = -19.9
$ws.Range("C16").Value = ""
$ws.Range("D16").Value = -14.1
$ws.Range("E16").Value = -5.3
$ws.Range("F16").Value = 0.70956

$ws.Range("A17").Value = "RM 88"
$ws.Range("B17").Value = -19.9
$ws.Range("C17").Value = 12.5
$ws.Range("D17").Value = -13.9
$ws.Range("E17").Value = -5.3
$ws.Range("F17").Value = 0.71011

$ws.Range("A18").Value = "RM 90"
$ws.Range("B18").Value = -19.6
$ws.Range("C18").Value = 14.4
$ws.Range("D18").Value = -13.1
$ws.Range("E18").Value = -5.4
$ws.Range("F18").Value = 0.71073

$ws.Range("A19").Value = "RM 95"
$ws.Range("B19").Value = -19.1
$ws.Range("C19").Value = 12.5
$ws.Range("D19").Value = -15.2
$ws.Range("E19").Value = -8.4
$ws.Range("F19").Value = 0.71076

$ws.Range("A20").Value = "RM 103"
$ws.Range("B20").Value = -19.5
$ws.Range("C20").Value = ""
$ws.Range("D20").Value = -15.3
$ws.Range("E20").Value = -6.9
$ws.Range("F20").Value = 0.7106

$ws.Range("A21").Value = "RM 116"
$ws.Range("B21").Value = -19.4
$ws.Range("C21").Value = 11.2
$ws.Range("D21").Value = -14.7
$ws.Range("E21").Value = ""
$ws.Range("F21").Value = 0.70981

$ws.Range("A22").Value = "RM 120"
$ws.Range("B22").Value = -19.6
$ws.Range("C22").Value = 11.5
$ws.Range("D22").Value = -15.2
$ws.Range("E22").Value = ""
$ws.Range("F22").Value = 0.70941

$ws.Range("A23").Value = "RM 125"
$ws.Range("B23").Value = -20.6
$ws.Range("C23").Value = 13.2
$ws.Range("D23").Value = -15.5
$ws.Range("E23").Value = ""
$ws.Range("F23").Value = 0.70931

$ws.Range("A24").Value = "RM 134"
$ws.Range("B24").Value = -19
$ws.Range("C24").Value = 12.5
$ws.Range("D24").Value = -14
$ws.Range("E24").Value = ""
$ws.Range("F24").Value = 0.71015

$ws.Range("A25").Value = "RM 135"
$ws.Range("B25").Value = -18.9
$ws.Range("C25").Value = 12.7
$ws.Range("D25").Value = -14.3
$ws.Range("E25").Value = ""
$ws.Range("F25").Value = 0.7108

$ws.Range("A26").Value = "RM 137"
$ws.Range("B26").Value = -19.5
$ws.Range("C26").Value = 12
$ws.Range("D26").Value = -14.9
$ws.Range("E26").Value = -8.5
$ws.Range("F26").Value = 0.70925

$ws.Range("A27").Value = "RM 138"
$ws.Range("B27").Value = -19.3
$ws.Range("C27").Value = 13.5
$ws.Range("D27").Value = -15.4
$ws.Range("E27").Value = -6.1
$ws.Range("F27").Value = 0.7092

$ws.Range("A28").Value = "RM 140"
$ws.Range("B28").Value = -19.5
$ws.Range("C28").Value = ""
$ws.Range("D28").Value = -13.9
$ws.Range("E28").Value = -7
$ws.Range("F28").Value = 0.70963

$ws.Range("A29").Value = "RM 142a"
$ws.Range("B29").Value = -17.7
$ws.Range("C29").Value = ""
$ws.Range("D29").Value = -13.9
$ws.Range("E29").Value = ""
$ws.Range("F29").Value = 0.70942

$ws.Range("A30").Value = "RM 142b"
$ws.Range("B30").Value = -19.5
$ws.Range("C30").Value = 12.2
$ws.Range("D30").Value = -13.8
$ws.Range("E30").Value = ""
$ws.Range("F30").Value = 0.71055

$ws.Range("A31").Value = "RM 145"
$ws.Range("B31").Value = -19.5
$ws.Range("C31").Value = 10.7
$ws.Range("D31").Value = ""
$ws.Range("E31").Value = -7.1
$ws.Range("F31").Value = 0.7116

$ws.Range("A32").Value = "RM 146"
$ws.Range("B32").Value = -19.4
$ws.Range("C32").Value = ""
$ws.Range("D32").Value = -15.2
$ws.Range("E32").Value = -7.8
$ws.Range("F32").Value = ""

$ws.Range("A33").Value = "RM 156"
$ws.Range("B33").Value = -19.4
$ws.Range("C33").Value = 12
$ws.Range("D33").Value = -15.1
$ws.Range("E33").Value = -9.5
$ws.Range("F33").Value = ""

$ws.Range("A34").Value = "RM 158"
$ws.Range("B34").Value = -19.6
$ws.Range("C34").Value = 10.3
$ws.Range("D34").Value = -14
$ws.Range("E34").Value = -7.6
$ws.Range("F34").Value = 0.70933

$ws.Range("A35").Value = "RM 159"
$ws.Range("B35").Value = -19.2
$ws.Range("C35").Value = 11.3
$ws.Range("D35").Value = -14.4
$ws.Range("E35").Value = -7.4
$ws.Range("F35").Value = ""

$ws.Range("A36").Value = "RM 165"
$ws.Range("B36").Value = -19.1
$ws.Range("C36").Value = 14.3
$ws.Range("D36").Value = -14.6
$ws.Range("E36").Value = -8.6
$ws.Range("F36").Value = 0.71087

$ws.Range("A37").Value = "RM 167"
$ws.Range("B37").Value = -19.8
$ws.Range("C37").Value = 12.1
$ws.Range("D37").Value = -14.3
$ws.Range("E37").Value = -7.1
$ws.Range("F37").Value = 0.71078

$ws.Range("A38").Value = "RM 170"
$ws.Range("B38").Value = -19.2
$ws.Range("C38").Value = 11.5
$ws.Range("D38").Value = -14.7
$ws.Range("E38").Value = -9.2
$ws.Range("F38").Value = ""

$ws.Range("A39").Value = "RM 173"
$ws.Range("B39").Value = -19.8
$ws.Range("C39").Value = 11.7
$ws.Range("D39").Value = ""
$ws.Range("E39").Value = -7.7
$ws.Range("F39").Value = 0.70929

$ws.Range("A40").Value = "RM 178"
$ws.Range("B40").Value = -19.5
$ws.Range("C40").Value = ""
$ws.Range("D40").Value = -14.4
$ws.Range("E40").Value = -7.9
$ws.Range("F40").Value = 0.71102

$ws.Range("A41").Value = "RM 186"
$ws.Range("B41").Value = -18.1
$ws.Range("C41").Value = 13.9
$ws.Range("D41").Value = -15.9
$ws.Range("E41").Value = -10.2
$ws.Range("F41").Value = ""

$ws.Range("A42").Value = "RM 193"
$ws.Range("B42").Value = -19
$ws.Range("C42").Value = 12
$ws.Range("D42").Value = -14.6
$ws.Range("E42").Value = -6.8
$ws.Range("F42").Value = 0.71115

$ws.Range("A43").Value = "RM 197"
$ws.Range("B43").Value = -18.9
$ws.Range("C43").Value = 13
$ws.Range("D43").Value = -14.1
$ws.Range("E43").Value = -8.8
$ws.Range("F43").Value = 0.71152

$ws.Range("A44").Value = "RM 202"
$ws.Range("B44").Value = -18.8
$ws.Range("C44").Value = 12.6
$ws.Range("D44").Value = -15.3
$ws.Range("E44").Value = -8.9
$ws.Range("F44").Value = 0.73857

$ws.Range("A45").Value = "RM 207"
$ws.Range("B45").Value = -19.7
$ws.Range("C45").Value = 11.7
$ws.Range("D45").Value = -14.9
$ws.Range("E45").Value = ""
$ws.Range("F45").Value = 0.70941

$ws.Range("A46").Value = "RM 208"
$ws.Range("B46").Value = -19.5
$ws.Range("C46").Value = ""
$ws.Range("D46").Value = -14.9
$ws.Range("E46").Value = ""
$ws.Range("F46").Value = 0.71128

$ws.Range("A47").Value = "RM 215"
$ws.Range("B47").Value = -19.8
$ws.Range("C47").Value = 12.9
$ws.Range("D47").Value = ""
$ws.Range("E47").Value = -6.6
$ws.Range("F47").Value = 0.71117

$ws.Range("A48").Value = "RM 219"
$ws.Range("B48").Value = -19.3
$ws.Range("C48").Value = 12.7
$ws.Range("D48").Value = -12.8
$ws.Range("E48").Value = -8.3
$ws.Range("F48").Value = 0.71045

$ws.Range("A49").Value = "RM 221"
$ws.Range("B49").Value = -19.9
$ws.Range("C49").Value = 11.5
$ws.Range("D49").Value = -13.7
$ws.Range("E49").Value = -7.2
$ws.Range("F49").Value = 0.70924

$ws.Range("A50").Value = "RM 232"
$ws.Range("B50").Value = -19.7
$ws.Range("C50").Value = 10.7
$ws.Range("D50").Value = -15.6
$ws.Range("E50").Value = -8.8
$ws.Range("F50").Value = 0.71064

$ws.Range("A51").Value = "RM 233"
$ws.Range("B51").Value = -20.5
$ws.Range("C51").Value = ""
$ws.Range("D51").Value = -14.5
$ws.Range("E51").Value = -7.7
$ws.Range("F51").Value = 0.71076

$ws.Range("A52").Value = "SC 5"
$ws.Range("B52").Value = -20.2
$ws.Range("C52").Value = ""
$ws.Range("D52").Value = ""
$ws.Range("E52").Value = -5
$ws.Range("F52").Value = 0.70948

$ws.Range("A53").Value = "SC 66"
$ws.Range("B53").Value = -20.3
$ws.Range("C53").Value = 10.5
$ws.Range("D53").Value = -12.9
$ws.Range("E53").Value = -5.7
$ws.Range("F53").Value = 0.71211

$ws.Range("A54").Value = "SC 92"
$ws.Range("B54").Value = -17.2
$ws.Range("C54").Value = 14.3
$ws.Range("D54").Value = -14
$ws.Range("E54").Value = -6.3
$ws.Range("F54").Value = 0.71157

$ws.Range("A55").Value = "SC 101"
$ws.Range("B55").Value = -20.4
$ws.Range("C55").Value = 10
$ws.Range("D55").Value = ""
$ws.Range("E55").Value = -10
$ws.Range("F55").Value = ""

$ws.Range("A56").Value = "SC 103"
$ws.Range("B56").Value = -19.2
$ws.Range("C56").Value = 11.9
$ws.Range("D56").Value = -14.7
$ws.Range("E56").Value = -5.7
$ws.Range("F56").Value = 0.7108

$ws.Range("A57").Value = "SC 105"
$ws.Range("B57").Value = -19.6
$ws.Range("C57").Value = ""
$ws.Range("D57").Value = -13.7
$ws.Range("E57").Value = -5.9
$ws.Range("F57").Value = 0.71146

$ws.Range("A58").Value = "SC 119"
$ws.Range("B58").Value = -19.5
$ws.Range("C58").Value = 11.2
$ws.Range("D58").Value = -13
$ws.Range("E58").Value = ""
$ws.Range("F58").Value = 0.70951

$ws.Range("A59").Value = "SC 120"
$ws.Range("B59").Value = -19.7
$ws.Range("C59").Value = ""
$ws.Range("D59").Value = -13.6
$ws.Range("E59").Value = -5.7
$ws.Range("F59").Value = 0.71096

$ws.Range("A60").Value = "SC 132"
$ws.Range("B60").Value = -18.8
$ws.Range("C60").Value = 15.3
$ws.Range("D60").Value = ""
$ws.Range("E60").Value = -8.1
$ws.Range("F60").Value = 0.70948

$ws.Range("A61").Value = "SC 193"
$ws.Range("B61").Value = -19.9
$ws.Range("C61").Value = 10.5
$ws.Range("D61").Value = ""
$ws.Range("E61").Value = -6.4
$ws.Range("F61").Value = 0.71183

$ws.Range("A62").Value = "SC 232"
$ws.Range("B62").Value = -19.5
$ws.Range("C62").Value = ""
$ws.Range("D62").Value = ""
$ws.Range("E62").Value = -10.7
$ws.Range("F62").Value = 0.71159

